$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false)
    if ($found) {
        $rng.Text = $new
    } else {
        Write-Output ("NOT FOUND: " + $old)
    }
}

Replace-Text "Umbizo limesahihishwa sio wakati" "Format has been corrected not the timing"
Replace-Text "Niliongeza sekunde 25 kwa kila muda ili kusahihisha wimbo wa utangulizi -john argentino" "I added 25 seconds to each timing to correct for the intro song -john argentino"
Replace-Text "Tatizo la uwanja wa ndege - manukuu:" "The airport problem - subtitles:"
Replace-Text "Utawala wa tatu" "The administrations of three"
Replace-Text "miji jirani: A, B na C waliamua" "neighboring cities: A, B and C decided"
Replace-Text "kujenga uwanja wa ndege unaogawanya gharama za" "to build an airport dividing the costs of"
Replace-Text "utekelezaji. Hali juu ya" "implementation. The condition on the"
Replace-Text "uchaguzi wa mahali pa kufaa zaidi ni" "choice of the most suitable place is"
Replace-Text "kwamba jumla ya umbali kutoka kwa kila mmoja" "that the sum of the distances from each"
Replace-Text "mji kwa uwanja wa ndege ni ndogo kama" "city to the airport is as small as"
Replace-Text "inawezekana. Timu ya wataalam wanaohusika" "possible. The team of experts in charge"
Replace-Text "ya kazi imeunda mfano wa kupata" "of the work has created a model to get"
Replace-Text "wazo la awali la mahali pa kuweka" "a preliminary idea of where to place the"
Replace-Text "muundo. Ovyo wao wapo" "structure. At their disposal there are"
Replace-Text "konokono wengine pete kubwa ya chuma na ndefu" "some snails a big metal ring and a long"
Replace-Text "kamba." "string."
Replace-Text "Eleza jinsi timu inaweza kusimamia matumizi" "Explain how the team can manage to use"
Replace-Text "nyenzo za kusema takriban" "the materials to tell approximately the"
Replace-Text "eneo bora la uwanja wa ndege. Fikiria" "ideal location of the airport. Imagine"
Replace-Text "kwamba miji imewekwa kwenye" "that the cities are placed at the"
Replace-Text "vipeo vya pembetatu ambayo ni" "vertices of a triangle which is"
Replace-Text "kwa hakika imetolewa tena kwa kiwango kama" "obviously reproduced in scale as"
Replace-Text "inavyoonyeshwa kwenye takwimu. Hili ni moja linalowezekana" "shown in figure. This is one possible"
Replace-Text "kuweka kamba huanza kutoka msumari mmoja," "setting the rope starts from one nail,"
Replace-Text "huenda ndani ya pete, huzunguka" "goes inside the ring, goes around the"
Replace-Text "msumari mwingine, msumari wa tatu, ndani ya" "other nail, the third nail, inside the"
Replace-Text "pete tena na sasa unaweza kuvuta tu" "ring again and now you can just pull the"
Replace-Text "kamba ili kupata uhakika huo" "rope in order to find the point that"
Replace-Text "unatafuta. Ili kufikia" "you're looking for. In order to reach the"
Replace-Text "uhakika, tunapaswa kusonga kamba kidogo" "point, we have to move the rope a bit"
Replace-Text "kwa sababu kuna " "because there is some "
Replace-Text "upinzani" "resistance"
Replace-Text " uliosababishwa" " caused"
Replace-Text "kwa nyenzo ambazo tunatumia lakini" "by the materials that we are using but"
Replace-Text "baada ya muda utafikia nafasi kutoka" "after a while you'll reach a position from"
Replace-Text "ambayo pete haisogei tena," "which the ring doesn't move anymore,"
Replace-Text "ambayo ni zaidi au chini ya hii. Na kama" "which is more or less this one. And as"
Replace-Text "kati ya pete na misumari ni" "between the ring and the nails are"
Replace-Text "kuwekwa zaidi au chini ya digrii 120 kutoka kwa moja" "placed more or less 120 degrees from one"
Replace-Text "nyingine ambayo ni 1/3 ya mduara," "another which is 1/3 of a circumference,"
Replace-Text "na hiyo ndiyo hatua tunayoiangalia" "and that's the point that we're looking"
Replace-Text "kwa: umbali wa chini kati ya" "for: the minimum distance between the"
Replace-Text "misumari na uwanja wa ndege unapojumlisha" "nails and the airport when you sum it"
Replace-Text "pamoja" "ogether"
Replace-Text "[Muziki]" "[Music]"
